{"js": "// Update the worksheet date heading and every two-digit multiplication\n// answer in the table (26 text replacements total), matching the target\n// commit. Cells are addressed by (row, col) position rather than by\n// searching for the old text, because a couple of old/new values collide\n// (e.g. \"72\u00d750=3600\" is both a prior answer and a new answer elsewhere),\n// which would make a naive global find/replace unsafe.\n\nconst body = context.document.body;\n\n// --- 1. Title paragraph: \"2023-11-06 Monday\" -> \"2023-11-07 Tuesday\" ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items[0].insertText(\"2023-11-07 Tuesday\", \"Replace\");\n\n// --- 2. Table answers ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Data rows (0-indexed) that contain answers; the rows in between are\n// blank spacer rows left untouched.\nconst updates = [\n  { row: 0, cells: [\"72\u00d750=3600\", \"77\u00d768=5236\", \"19\u00d785=1615\", \"46\u00d722=1012\", \"60\u00d734=2040\"] },\n  { row: 4, cells: [\"20\u00d726=520\", \"38\u00d798=3724\", \"16\u00d748=768\", \"54\u00d747=2538\", \"67\u00d780=5360\"] },\n  { row: 9, cells: [\"52\u00d779=4108\", \"80\u00d745=3600\", \"29\u00d758=1682\", \"61\u00d765=3965\", \"54\u00d778=4212\"] },\n  { row: 14, cells: [\"51\u00d733=1683\", \"15\u00d727=405\", \"53\u00d765=3445\", \"44\u00d725=1100\", \"70\u00d784=5880\"] },\n  { row: 19, cells: [\"46\u00d799=4554\", \"33\u00d728=924\", \"34\u00d762=2108\", \"81\u00d798=7938\", \"67\u00d784=5628\"] },\n];\n\nfor (const { row, cells } of updates) {\n  for (let col = 0; col < cells.length; col++) {\n    const cell = table.getCell(row, col);\n    cell.body.paragraphs.load(\"items\");\n  }\n}\nawait context.sync();\n\nfor (const { row, cells } of updates) {\n  for (let col = 0; col < cells.length; col++) {\n    const cell = table.getCell(row, col);\n    cell.body.paragraphs.items[0].insertText(cells[col], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date heading and every two-digit multiplication\n# answer in the table (26 text replacements total), matching the target\n# commit. Table cells are addressed by their (row, col) position rather\n# than by searching for the old text, because a couple of old/new values\n# collide (e.g. \"72x50=3600\" is both a prior answer and a new answer\n# elsewhere), which would make a naive global find/replace unsafe.\n\n$d = $word.ActiveDocument\n\n# --- 1. Title paragraph: \"2023-11-06 Monday\" -> \"2023-11-07 Tuesday\" ---\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Execute(\"2023-11-06 Monday\", $false, $false, $false, $false, $false, $true, 1, $false, \"2023-11-07 Tuesday\", 2) | Out-Null\n\n# --- 2. Table answers, addressed by (row, col), 1-indexed per COM OM ---\n$tbl = $d.Tables(1)\n\n$updates = @(\n    @{ Row = 1;  Values = @(\"72\u00d750=3600\", \"77\u00d768=5236\", \"19\u00d785=1615\", \"46\u00d722=1012\", \"60\u00d734=2040\") },\n    @{ Row = 5;  Values = @(\"20\u00d726=520\",  \"38\u00d798=3724\", \"16\u00d748=768\",  \"54\u00d747=2538\", \"67\u00d780=5360\") },\n    @{ Row = 10; Values = @(\"52\u00d779=4108\", \"80\u00d745=3600\", \"29\u00d758=1682\", \"61\u00d765=3965\", \"54\u00d778=4212\") },\n    @{ Row = 15; Values = @(\"51\u00d733=1683\", \"15\u00d727=405\",  \"53\u00d765=3445\", \"44\u00d725=1100\", \"70\u00d784=5880\") },\n    @{ Row = 20; Values = @(\"46\u00d799=4554\", \"33\u00d728=924\",  \"34\u00d762=2108\", \"81\u00d798=7938\", \"67\u00d784=5628\") }\n)\n\nforeach ($update in $updates) {\n    $row = $update.Row\n    $values = $update.Values\n    for ($col = 1; $col -le $values.Count; $col++) {\n        $tbl.Cell($row, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
